# Fruta / hortaliza, semanal
#
# The weekly refresh re-shuffles the "Fecha" (D) and the data-series
# columns (M..T: Volumen, Precio minimo/maximo/promedio, Unidad de
# comercializacion, Origen, Precio $/Kg, Kg/unidad) across the existing
# data rows (2..26). The identifying columns A,B,C,E..L (Mercado,
# Region, Codreg, Tipo, Producto, Categoria, Variedad, Calidad) stay put.
#
# Row r's new D..T-ish payload == the old payload that used to live on
# row $sigma[r] (a permutation derived from the published diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sigma = @{
  2  = 18
  3  = 20
  4  = 5
  5  = 16
  6  = 6
  7  = 8
  8  = 22
  9  = 14
  10 = 23
  11 = 9
  12 = 3
  13 = 13
  14 = 2
  15 = 11
  16 = 25
  17 = 12
  18 = 10
  19 = 26
  20 = 17
  21 = 15
  22 = 19
  23 = 24
  24 = 4
  25 = 7
  26 = 21
}

# Columns that move with the row: D(4), M..T(13..20)
$cols = @(4, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot every source value first (Value2 avoids COM date/variant
# wrapping) so writes never clobber data still needed for a later read.
$snapshot = @{}
for ($r = 2; $r -le 26; $r++) {
  $rowVals = @{}
  foreach ($c in $cols) {
    $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
  }
  $snapshot[$r] = $rowVals
}

for ($r = 2; $r -le 26; $r++) {
  $src = $sigma[$r]
  $rowVals = $snapshot[$src]
  foreach ($c in $cols) {
    $ws.Cells.Item($r, $c).Value2 = $rowVals[$c]
  }
}
